$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the three new worksheets (GBIFTaxa, NCBITaxa, Locations) after the
#    existing "Summary" sheet, in that order.
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsGBIF = $wb.Worksheets.Add($null, $afterSheet)
$wsGBIF.Name = "GBIFTaxa"

$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNCBI = $wb.Worksheets.Add($null, $afterSheet)
$wsNCBI.Name = "NCBITaxa"

$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLoc = $wb.Worksheets.Add($null, $afterSheet)
$wsLoc.Name = "Locations"

$wsSummary = $wb.Worksheets.Item("Summary")

# ---------------------------------------------------------------------------
# 2. Populate GBIFTaxa
# ---------------------------------------------------------------------------
$wsGBIF.Range("A1").Value = "Name"
$wsGBIF.Range("B1").Value = "Taxon name"
$wsGBIF.Range("C1").Value = "Taxon type"
$wsGBIF.Range("D1").Value = "Taxon ID"
$wsGBIF.Range("E1").Value = "Ignore ID"
$wsGBIF.Range("F1").Value = "Parent name"
$wsGBIF.Range("G1").Value = "Parent type"
$wsGBIF.Range("H1").Value = "Parent ID"
$wsGBIF.Range("I1").Value = "Comments"

$wsGBIF.Range("A2").Value = "crbo"
$wsGBIF.Range("B2").Value = "Crematogaster borneensis"
$wsGBIF.Range("C2").Value = "Species"

$wsGBIF.Range("A3").Value = "dolic_sp"
$wsGBIF.Range("B3").Value = "Dolichoderus"
$wsGBIF.Range("C3").Value = "Genus"

$wsGBIF.Range("A4").Value = "gannets"
$wsGBIF.Range("B4").Value = "Morus"
$wsGBIF.Range("C4").Value = "Genus"
$wsGBIF.Range("D4").Value = 2480962

$wsGBIF.Range("A5").Value = "lost_orang"
$wsGBIF.Range("B5").Value = "Pongo tapanuliensis"
$wsGBIF.Range("C5").Value = "Species"
$wsGBIF.Range("F5").Value = "Pongo"
$wsGBIF.Range("G5").Value = "Genus"
$wsGBIF.Range("I5").Value = "New species"

$wsGBIF.Range("A6").Value = "morpho1"
$wsGBIF.Range("B6").Value = "NA"
$wsGBIF.Range("C6").Value = "Morphospecies"
$wsGBIF.Range("F6").Value = "Formicidae"
$wsGBIF.Range("G6").Value = "Family"

$wsGBIF.Range("A7").Value = "bombines"
$wsGBIF.Range("B7").Value = "Bombini"
$wsGBIF.Range("C7").Value = "Tribe"
$wsGBIF.Range("F7").Value = "Apidae"
$wsGBIF.Range("G7").Value = "Family"

$wsGBIF.Range("A8").Value = "micr_hid"
$wsGBIF.Range("B8").Value = "Microcopris hidakai"
$wsGBIF.Range("C8").Value = "Species"
$wsGBIF.Range("E8").Value = 1090433
$wsGBIF.Range("F8").Value = "Microcopris"
$wsGBIF.Range("G8").Value = "Genus"

# ---------------------------------------------------------------------------
# 3. Populate NCBITaxa
# ---------------------------------------------------------------------------
$wsNCBI.Range("A1").Value = "Name"
$wsNCBI.Range("B1").Value = "New"
$wsNCBI.Range("C1").Value = "Superkingdom"
$wsNCBI.Range("D1").Value = "Kingdom"
$wsNCBI.Range("E1").Value = "Phylum"
$wsNCBI.Range("F1").Value = "Class"
$wsNCBI.Range("G1").Value = "Comments"

$wsNCBI.Range("A2").Value = "G_proteobacteria"
$wsNCBI.Range("C2").Value = "Bacteria"
$wsNCBI.Range("E2").Value = "Pseudomonadota"
$wsNCBI.Range("F2").Value = "Gammaproteobacteria"

$wsNCBI.Range("A3").Value = "E_mycetes"
$wsNCBI.Range("C3").Value = "Eukaryota"
$wsNCBI.Range("D3").Value = "Fungi"
$wsNCBI.Range("E3").Value = "Ascomycota"
$wsNCBI.Range("F3").Value = "Eurotiomycetes"

$wsNCBI.Range("A4").Value = "Dinophyceae"
$wsNCBI.Range("C4").Value = "Eukaryota"
$wsNCBI.Range("F4").Value = "Dinophyceae"

$wsNCBI.Range("A5").Value = "Acidobact"
$wsNCBI.Range("D5").Value = "k__Bacteria"
$wsNCBI.Range("E5").Value = "p__Acidobacteria"
$wsNCBI.Range("F5").Value = "c__Acidobacteriia"

$wsNCBI.Range("A6").Value = "New_fungus"
$wsNCBI.Range("B6").Value = "Yes"
$wsNCBI.Range("C6").Value = "Eukaryota"
$wsNCBI.Range("D6").Value = "Fungi"
$wsNCBI.Range("E6").Value = "Ascomycota"
$wsNCBI.Range("F6").Value = "Mynewfungusetes"

# ---------------------------------------------------------------------------
# 4. Populate Locations
# ---------------------------------------------------------------------------
$wsLoc.Range("A1").Value = "Location name"
$wsLoc.Range("B1").Value = "New"
$wsLoc.Range("C1").Value = "Latitude"
$wsLoc.Range("D1").Value = "Longitude"
$wsLoc.Range("E1").Value = "Type"
$wsLoc.Range("F1").Value = "WKT"

$wsLoc.Range("A2").Value = "E_194"
$wsLoc.Range("B2").Value = "No"

$wsLoc.Range("A3").Value = "E_195"
$wsLoc.Range("B3").Value = "No"

$wsLoc.Range("A4").Value = "My_site_1"
$wsLoc.Range("B4").Value = "Yes"
$wsLoc.Range("C4").Value = 4.9577210000000003
$wsLoc.Range("D4").Value = 117.776023
$wsLoc.Range("E4").Value = "POINT"
$wsLoc.Range("F4").Value = "NA"

$wsLoc.Range("A5").Value = "My_site_2"
$wsLoc.Range("B5").Value = "Yes"
$wsLoc.Range("C5").Value = "NA"
$wsLoc.Range("D5").Value = "NA"
$wsLoc.Range("E5").Value = "POINT"
$wsLoc.Range("F5").Value = "NA"

$wsLoc.Range("A6").Value = "My_site_3"
$wsLoc.Range("B6").Value = "Yes"
$wsLoc.Range("C6").Value = "NA"
$wsLoc.Range("D6").Value = "NA"
$wsLoc.Range("E6").Value = "POINT"
$wsLoc.Range("F6").Value = "Point(117.7762 4.9576)"

$wsLoc.Range("A7").Value = "My_transect_1"
$wsLoc.Range("B7").Value = "Yes"
$wsLoc.Range("C7").Value = "NA"
$wsLoc.Range("D7").Value = "NA"
$wsLoc.Range("E7").Value = "Linestring"
$wsLoc.Range("F7").Value = "Linestring(117.7762 4.9576, 117.7862 4.9676)"

# ---------------------------------------------------------------------------
# 5. Column widths (best-fit style autosizing approximated via ColumnWidth)
# ---------------------------------------------------------------------------
$wsSummary.Columns.Item(2).ColumnWidth = 41.666666666666664   # B -> 42.5

$wsGBIF.Columns.Item(2).ColumnWidth = 24.833333333333336      # B -> 25.6640625
$wsGBIF.Columns.Item(3).ColumnWidth = 14.333333333333332      # C -> 15.1640625
$wsGBIF.Columns.Item(4).ColumnWidth = 8.5                     # D -> 9.33203125
$wsGBIF.Columns.Item(6).ColumnWidth = 10.833333333333332      # F -> 11.6640625
$wsGBIF.Columns.Item(7).ColumnWidth = 9.666666666666666       # G -> 10.5

$wsNCBI.Columns.Item(1).ColumnWidth = 16.333333333333336      # A -> 17.1640625
$wsNCBI.Columns.Item(3).ColumnWidth = 12.0                    # C -> 12.83203125
$wsNCBI.Columns.Item(4).ColumnWidth = 11.166666666666666      # D -> 12
$wsNCBI.Columns.Item(5).ColumnWidth = 16.666666666666668      # E -> 17.5
$wsNCBI.Columns.Item(6).ColumnWidth = 21.0                    # F -> 21.83203125

$wsLoc.Columns.Item(1).ColumnWidth = 12.333333333333332       # A -> 13.1640625
$wsLoc.Columns.Item(2).ColumnWidth = 4.166666666666667        # B -> 5
$wsLoc.Columns.Item(3).ColumnWidth = 8.333333333333332        # C -> 9.1640625
$wsLoc.Columns.Item(4).ColumnWidth = 10.333333333333332       # D -> 11.1640625
$wsLoc.Columns.Item(5).ColumnWidth = 8.333333333333332        # E -> 9.1640625
$wsLoc.Columns.Item(6).ColumnWidth = 38.5                     # F -> 39.33203125

# ---------------------------------------------------------------------------
# 6. Sheet view / selection state.
#    Locations becomes the active/visible tab; Summary loses its selection.
# ---------------------------------------------------------------------------
$wsGBIF.Range("C5").Select()
$wsNCBI.Range("F17").Select()
$wsSummary.Range("A15").Select()
$wsLoc.Range("G27").Select()

$wsLoc.Activate()

$excel.ActiveWindow.ScrollColumn = 1
